{"js": "/*\n * Reproduce: add the three SharePoint/Content-Organizer custom XML parts\n * that Word's document-management integration stamps into the package\n * (customXml/item1.xml .. item3.xml + their itemProps companions).\n *\n * The only content-level surface Word exposes for this is\n * `document.customXmlParts.add(xml)` \u2014 each call mints one new\n * Custom XML Part (and Word itself writes out the matching\n * customXml/itemPropsN.xml + relationships + content-type entries\n * when the package is serialized). No existing body/style content is\n * touched by this change, matching the source diff which only adds\n * these three parts.\n */\n\nconst customXmlParts = context.document.customXmlParts;\n\n// 1) Content-type / managed-metadata schema part (SharePoint \"Document\" content type).\nconst contentTypeSchemaXml = \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<ct:contentTypeSchema xmlns:ct=\\\"http://schemas.microsoft.com/office/2006/metadata/contentType\\\" xmlns:ma=\\\"http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes\\\" ct:_=\\\"\\\" ma:_=\\\"\\\" ma:contentTypeName=\\\"Document\\\" ma:contentTypeID=\\\"0x010100ACF816946BF9B24DB7F26955C4543618\\\" ma:contentTypeVersion=\\\"9\\\" ma:contentTypeDescription=\\\"Create a new document.\\\" ma:contentTypeScope=\\\"\\\" ma:versionID=\\\"f83a017b00179535ba734d13d836873f\\\">\\n  <xsd:schema xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:p=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" xmlns:ns2=\\\"4526b749-01de-4416-a7d8-a2effcd30547\\\" targetNamespace=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" ma:root=\\\"true\\\" ma:fieldsID=\\\"2f3f2e2b342ee0fc112d7f8f6bff53ff\\\" ns2:_=\\\"\\\">\\n    <xsd:import namespace=\\\"4526b749-01de-4416-a7d8-a2effcd30547\\\"/>\\n    <xsd:element name=\\\"properties\\\">\\n      <xsd:complexType>\\n        <xsd:sequence>\\n          <xsd:element name=\\\"documentManagement\\\">\\n            <xsd:complexType>\\n              <xsd:all>\\n                <xsd:element ref=\\\"ns2:MediaServiceMetadata\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceFastMetadata\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceAutoKeyPoints\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceKeyPoints\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceAutoTags\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceOCR\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceGenerationTime\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceEventHashCode\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceDateTaken\\\" minOccurs=\\\"0\\\"/>\\n              </xsd:all>\\n            </xsd:complexType>\\n          </xsd:element>\\n        </xsd:sequence>\\n      </xsd:complexType>\\n    </xsd:element>\\n  </xsd:schema>\\n  <xsd:schema xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:dms=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\" xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" targetNamespace=\\\"4526b749-01de-4416-a7d8-a2effcd30547\\\" elementFormDefault=\\\"qualified\\\">\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\"/>\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\"/>\\n    <xsd:element name=\\\"MediaServiceMetadata\\\" ma:index=\\\"8\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceMetadata\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceMetadata\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceFastMetadata\\\" ma:index=\\\"9\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceFastMetadata\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceFastMetadata\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceAutoKeyPoints\\\" ma:index=\\\"10\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceAutoKeyPoints\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceAutoKeyPoints\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceKeyPoints\\\" ma:index=\\\"11\\\" nillable=\\\"true\\\" ma:displayName=\\\"KeyPoints\\\" ma:internalName=\\\"MediaServiceKeyPoints\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\">\\n          <xsd:maxLength value=\\\"255\\\"/>\\n        </xsd:restriction>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceAutoTags\\\" ma:index=\\\"12\\\" nillable=\\\"true\\\" ma:displayName=\\\"Tags\\\" ma:internalName=\\\"MediaServiceAutoTags\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceOCR\\\" ma:index=\\\"13\\\" nillable=\\\"true\\\" ma:displayName=\\\"Extracted Text\\\" ma:internalName=\\\"MediaServiceOCR\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\">\\n          <xsd:maxLength value=\\\"255\\\"/>\\n        </xsd:restriction>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceGenerationTime\\\" ma:index=\\\"14\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceGenerationTime\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceGenerationTime\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceEventHashCode\\\" ma:index=\\\"15\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceEventHashCode\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceEventHashCode\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceDateTaken\\\" ma:index=\\\"16\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceDateTaken\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceDateTaken\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n  </xsd:schema>\\n  <xsd:schema xmlns=\\\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\\\" xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xsi=\\\"http://www.w3.org/2001/XMLSchema-instance\\\" xmlns:dc=\\\"http://purl.org/dc/elements/1.1/\\\" xmlns:dcterms=\\\"http://purl.org/dc/terms/\\\" xmlns:odoc=\\\"http://schemas.microsoft.com/internal/obd\\\" targetNamespace=\\\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\\\" elementFormDefault=\\\"qualified\\\" attributeFormDefault=\\\"unqualified\\\" blockDefault=\\\"#all\\\">\\n    <xsd:import namespace=\\\"http://purl.org/dc/elements/1.1/\\\" schemaLocation=\\\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd\\\"/>\\n    <xsd:import namespace=\\\"http://purl.org/dc/terms/\\\" schemaLocation=\\\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd\\\"/>\\n    <xsd:element name=\\\"coreProperties\\\" type=\\\"CT_coreProperties\\\"/>\\n    <xsd:complexType name=\\\"CT_coreProperties\\\">\\n      <xsd:all>\\n        <xsd:element ref=\\\"dc:creator\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dcterms:created\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dc:identifier\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"contentType\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\" ma:index=\\\"0\\\" ma:displayName=\\\"Content Type\\\"/>\\n        <xsd:element ref=\\\"dc:title\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" ma:index=\\\"4\\\" ma:displayName=\\\"Title\\\"/>\\n        <xsd:element ref=\\\"dc:subject\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dc:description\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"keywords\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element ref=\\\"dc:language\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"category\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element name=\\\"version\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element name=\\\"revision\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\">\\n          <xsd:annotation>\\n            <xsd:documentation>\\n                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.\\n                    </xsd:documentation>\\n          </xsd:annotation>\\n        </xsd:element>\\n        <xsd:element name=\\\"lastModifiedBy\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element ref=\\\"dcterms:modified\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"contentStatus\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n      </xsd:all>\\n    </xsd:complexType>\\n  </xsd:schema>\\n  <xs:schema xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" targetNamespace=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" elementFormDefault=\\\"qualified\\\" attributeFormDefault=\\\"unqualified\\\">\\n    <xs:element name=\\\"Person\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:DisplayName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:AccountId\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:AccountType\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"DisplayName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"AccountId\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"AccountType\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"BDCAssociatedEntity\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:BDCEntity\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"unbounded\\\"/>\\n        </xs:sequence>\\n        <xs:attribute ref=\\\"pc:EntityNamespace\\\"/>\\n        <xs:attribute ref=\\\"pc:EntityName\\\"/>\\n        <xs:attribute ref=\\\"pc:SystemInstanceName\\\"/>\\n        <xs:attribute ref=\\\"pc:AssociationName\\\"/>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:attribute name=\\\"EntityNamespace\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"EntityName\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"SystemInstanceName\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"AssociationName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"BDCEntity\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:EntityDisplayName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityInstanceReference\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId1\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId2\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId3\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId4\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId5\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"EntityDisplayName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityInstanceReference\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId1\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId2\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId3\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId4\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId5\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"Terms\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:TermInfo\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"unbounded\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"TermInfo\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:TermName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:TermId\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"TermName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"TermId\\\" type=\\\"xs:string\\\"/>\\n  </xs:schema>\\n</ct:contentTypeSchema>\\n\";\n\n// 2) SharePoint list-form template part.\nconst formTemplatesXml = \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<?mso-contentType ?>\\n<FormTemplates xmlns=\\\"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms\\\">\\n  <Display>DocumentLibraryForm</Display>\\n  <Edit>DocumentLibraryForm</Edit>\\n  <New>DocumentLibraryForm</New>\\n</FormTemplates>\\n\";\n\n// 3) Empty documentManagement properties part.\nconst propertiesXml = \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<p:properties xmlns:p=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" xmlns:xsi=\\\"http://www.w3.org/2001/XMLSchema-instance\\\" xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\">\\n  <documentManagement/>\\n</p:properties>\\n\";\n\ncustomXmlParts.add(contentTypeSchemaXml);\ncustomXmlParts.add(formTemplatesXml);\ncustomXmlParts.add(propertiesXml);\n\nawait context.sync();\n", "ps1": "# Reproduce: add the three SharePoint/Content-Organizer custom XML parts\n# that Word's document-management integration stamps into the package\n# (customXml/item1.xml .. item3.xml + their itemProps companions).\n#\n# The only content-level surface Word automation exposes for this is\n# Document.CustomXMLParts.Add(xml) -- each call mints one new Custom XML\n# Part (Word itself writes out the matching customXml/itemPropsN.xml,\n# relationships and content-type entries when the package is saved).\n# No existing body/style content is touched by this change, matching the\n# source diff which only adds these three parts.\n\n$d = $word.ActiveDocument\n$parts = $d.CustomXMLParts\n\n# 1) Content-type / managed-metadata schema part (SharePoint \"Document\" content type).\n$contentTypeSchemaXml = '<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<ct:contentTypeSchema xmlns:ct=\"http://schemas.microsoft.com/office/2006/metadata/contentType\" xmlns:ma=\"http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes\" ct:_=\"\" ma:_=\"\" ma:contentTypeName=\"Document\" ma:contentTypeID=\"0x010100ACF816946BF9B24DB7F26955C4543618\" ma:contentTypeVersion=\"9\" ma:contentTypeDescription=\"Create a new document.\" ma:contentTypeScope=\"\" ma:versionID=\"f83a017b00179535ba734d13d836873f\">\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:ns2=\"4526b749-01de-4416-a7d8-a2effcd30547\" targetNamespace=\"http://schemas.microsoft.com/office/2006/metadata/properties\" ma:root=\"true\" ma:fieldsID=\"2f3f2e2b342ee0fc112d7f8f6bff53ff\" ns2:_=\"\">\n    <xsd:import namespace=\"4526b749-01de-4416-a7d8-a2effcd30547\"/>\n    <xsd:element name=\"properties\">\n      <xsd:complexType>\n        <xsd:sequence>\n          <xsd:element name=\"documentManagement\">\n            <xsd:complexType>\n              <xsd:all>\n                <xsd:element ref=\"ns2:MediaServiceMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceFastMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceAutoKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceAutoTags\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceOCR\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceGenerationTime\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceEventHashCode\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceDateTaken\" minOccurs=\"0\"/>\n              </xsd:all>\n            </xsd:complexType>\n          </xsd:element>\n        </xsd:sequence>\n      </xsd:complexType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"4526b749-01de-4416-a7d8-a2effcd30547\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"MediaServiceMetadata\" ma:index=\"8\" nillable=\"true\" ma:displayName=\"MediaServiceMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceFastMetadata\" ma:index=\"9\" nillable=\"true\" ma:displayName=\"MediaServiceFastMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceFastMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceAutoKeyPoints\" ma:index=\"10\" nillable=\"true\" ma:displayName=\"MediaServiceAutoKeyPoints\" ma:hidden=\"true\" ma:internalName=\"MediaServiceAutoKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceKeyPoints\" ma:index=\"11\" nillable=\"true\" ma:displayName=\"KeyPoints\" ma:internalName=\"MediaServiceKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceAutoTags\" ma:index=\"12\" nillable=\"true\" ma:displayName=\"Tags\" ma:internalName=\"MediaServiceAutoTags\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceOCR\" ma:index=\"13\" nillable=\"true\" ma:displayName=\"Extracted Text\" ma:internalName=\"MediaServiceOCR\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceGenerationTime\" ma:index=\"14\" nillable=\"true\" ma:displayName=\"MediaServiceGenerationTime\" ma:hidden=\"true\" ma:internalName=\"MediaServiceGenerationTime\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceEventHashCode\" ma:index=\"15\" nillable=\"true\" ma:displayName=\"MediaServiceEventHashCode\" ma:hidden=\"true\" ma:internalName=\"MediaServiceEventHashCode\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceDateTaken\" ma:index=\"16\" nillable=\"true\" ma:displayName=\"MediaServiceDateTaken\" ma:hidden=\"true\" ma:internalName=\"MediaServiceDateTaken\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:dc=\"http://purl.org/dc/elements/1.1/\" xmlns:dcterms=\"http://purl.org/dc/terms/\" xmlns:odoc=\"http://schemas.microsoft.com/internal/obd\" targetNamespace=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\" blockDefault=\"#all\">\n    <xsd:import namespace=\"http://purl.org/dc/elements/1.1/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd\"/>\n    <xsd:import namespace=\"http://purl.org/dc/terms/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd\"/>\n    <xsd:element name=\"coreProperties\" type=\"CT_coreProperties\"/>\n    <xsd:complexType name=\"CT_coreProperties\">\n      <xsd:all>\n        <xsd:element ref=\"dc:creator\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dcterms:created\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:identifier\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentType\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\" ma:index=\"0\" ma:displayName=\"Content Type\"/>\n        <xsd:element ref=\"dc:title\" minOccurs=\"0\" maxOccurs=\"1\" ma:index=\"4\" ma:displayName=\"Title\"/>\n        <xsd:element ref=\"dc:subject\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:description\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"keywords\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dc:language\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"category\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"version\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"revision\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\">\n          <xsd:annotation>\n            <xsd:documentation>\n                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.\n                    </xsd:documentation>\n          </xsd:annotation>\n        </xsd:element>\n        <xsd:element name=\"lastModifiedBy\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dcterms:modified\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentStatus\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n      </xsd:all>\n    </xsd:complexType>\n  </xsd:schema>\n  <xs:schema xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" targetNamespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\">\n    <xs:element name=\"Person\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:DisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountId\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountType\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"DisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"AccountId\" type=\"xs:string\"/>\n    <xs:element name=\"AccountType\" type=\"xs:string\"/>\n    <xs:element name=\"BDCAssociatedEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:BDCEntity\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n        <xs:attribute ref=\"pc:EntityNamespace\"/>\n        <xs:attribute ref=\"pc:EntityName\"/>\n        <xs:attribute ref=\"pc:SystemInstanceName\"/>\n        <xs:attribute ref=\"pc:AssociationName\"/>\n      </xs:complexType>\n    </xs:element>\n    <xs:attribute name=\"EntityNamespace\" type=\"xs:string\"/>\n    <xs:attribute name=\"EntityName\" type=\"xs:string\"/>\n    <xs:attribute name=\"SystemInstanceName\" type=\"xs:string\"/>\n    <xs:attribute name=\"AssociationName\" type=\"xs:string\"/>\n    <xs:element name=\"BDCEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:EntityDisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityInstanceReference\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId1\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId2\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId3\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId4\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId5\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"EntityDisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"EntityInstanceReference\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId1\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId2\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId3\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId4\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId5\" type=\"xs:string\"/>\n    <xs:element name=\"Terms\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermInfo\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:TermId\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermName\" type=\"xs:string\"/>\n    <xs:element name=\"TermId\" type=\"xs:string\"/>\n  </xs:schema>\n</ct:contentTypeSchema>\n'\n\n# 2) SharePoint list-form template part.\n$formTemplatesXml = '<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<?mso-contentType ?>\n<FormTemplates xmlns=\"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms\">\n  <Display>DocumentLibraryForm</Display>\n  <Edit>DocumentLibraryForm</Edit>\n  <New>DocumentLibraryForm</New>\n</FormTemplates>\n'\n\n# 3) Empty documentManagement properties part.\n$propertiesXml = '<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<p:properties xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\">\n  <documentManagement/>\n</p:properties>\n'\n\n$parts.Add($contentTypeSchemaXml) | Out-Null\n$parts.Add($formTemplatesXml) | Out-Null\n$parts.Add($propertiesXml) | Out-Null\n"}
